$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the "Real time (minutes)" values for the Transfer-related tasks
$ws.Range("D30").Value = 15
$ws.Range("D31").Value = 45
$ws.Range("D32").Value = 5
$ws.Range("D33").Value = 10

# Update the active selection on the sheet (as left by the editor)
$ws.Range("D40").Select()
